{"js": "// Load the paragraphs of the document body.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\n// Paragraph layout before the edit:\n//   0: \"Write Up\"   (Title style)\n//   1: \"\"            (empty, default/Normal style)\n//   2: \"\"            (empty, Heading1 style) -> becomes the section title\nconst blankParagraph = items[1];\nconst headingParagraph = items[items.length - 1];\n\n// Insert the two new intro paragraphs right after the blank paragraph\n// (and therefore right before the heading paragraph), so they pick up\n// the default/Normal style instead of inheriting Heading1.\nlet insertAfter = blankParagraph;\n\nconst introParagraph = insertAfter.insertParagraph(\n  \"In this tutorial, we will be throwing a door into the town, so when we are kicked out, when the hero dies, we can get back in without having to restart the entire game. We will also be writing a small bit of code, to create a keypress event to get back to the fight room. The fight room is not using a door, but a monster collision event, so to test it we will simply create a temporary F1 keypress.\",\n  Word.InsertLocation.after\n);\ninsertAfter = introParagraph;\n\ninsertAfter.insertParagraph(\n  \"So, if this sounds interesting to you, then please join us for our brand-new tutorial entitled:\",\n  Word.InsertLocation.after\n);\n\n// Fill in the text of the heading paragraph itself.\nheadingParagraph.insertText(\"19 Getting Back in the Door\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$introText = \"In this tutorial, we will be throwing a door into the town, so when we are kicked out, when the hero dies, we can get back in without having to restart the entire game. We will also be writing a small bit of code, to create a keypress event to get back to the fight room. The fight room is not using a door, but a monster collision event, so to test it we will simply create a temporary F1 keypress.\"\n$calloutText = \"So, if this sounds interesting to you, then please join us for our brand-new tutorial entitled:\"\n$headingText = \"19 Getting Back in the Door\"\n\n# Paragraph 1 = \"Write Up\" (Title), Paragraph 2 = blank (Normal),\n# Paragraph 3 (last) = blank (Heading 1) before this edit runs.\n\n# Insert the two new introductory paragraphs right after the blank\n# \"Normal\" styled paragraph, so each one picks up the default style\n# instead of inheriting Heading1 from the paragraph that follows them.\n$blank = $d.Paragraphs(2)\n$blank.Range.InsertParagraphAfter()\n$d.Paragraphs(3).Range.Text = $introText\n\n$d.Paragraphs(3).Range.InsertParagraphAfter()\n$d.Paragraphs(4).Range.Text = $calloutText\n\n# Fill in the text of the (still) last paragraph, the Heading1 title.\n$headingParagraph = $d.Paragraphs($d.Paragraphs.Count)\n$headingParagraph.Range.Text = $headingText\n"}
